# Senior-project weekly report update: refresh the Time log with the
# latest work entries, update the Report sheet's rolled-up totals and
# narrative text (Accomplishments / Problems / Scripture / etc.), make
# "Time" the active/selected sheet, and drop the now-stale log_2 defined
# name that pointed at the external log.csv query.

$wb = $excel.ActiveWorkbook
$report = $wb.Worksheets.Item("Report")
$time = $wb.Worksheets.Item("Time")

# ---------------------------------------------------------------------
# 1. Time sheet: append the new log rows (95-106)
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row = 95;  Task = "Math Unit Test and code";                Start = 41730.887291666666; Stop = 41730.892222222225 },
    @{ Row = 96;  Task = "Sequential Ray Tracer Tests and code";    Start = 41730.892222222225; Stop = 41730.908935185187 },
    @{ Row = 97;  Task = "Design";                                  Start = 41731.749699074076; Stop = 41731.751331018517 },
    @{ Row = 98;  Task = "Research";                                Start = 41731.751331018517; Stop = 41731.755520833336 },
    @{ Row = 99;  Task = "Sequential Ray Tracer Tests and code";    Start = 41731.755532407406; Stop = 41731.778969907406 },
    @{ Row = 100; Task = "Research";                                Start = 41731.778969907406; Stop = 41731.804942129631 },
    @{ Row = 101; Task = "Design";                                  Start = 41731.804942129631; Stop = 41731.811168981483 },
    @{ Row = 102; Task = "Sequential Ray Tracer Tests and code";    Start = 41731.811168981483; Stop = 41731.828935185185 },
    @{ Row = 103; Task = "Sequential Ray Tracer Tests and code";    Start = 41731.832708333335; Stop = 41731.847071759257 },
    @{ Row = 104; Task = "Sequential Ray Tracer Tests and code";    Start = 41731.867361111108; Stop = 41731.892604166664 },
    @{ Row = 105; Task = "Design";                                  Start = 41731.892604166664; Stop = 41731.897997685184 },
    @{ Row = 106; Task = "Research";                                Start = 41731.898009259261; Stop = 41731.9143287037 }
)

# Row 94's Stop time (C94) used to be a live =NOW() snapshot; freeze it
# to the value it held at the moment this report was written, and carry
# its formatting down into the freshly appended rows.
$time.Range("C94").Formula = $null
$time.Range("C94").Value = 41730.887291666666
$time.Range("B94:D94").Copy() | Out-Null
$time.Range("B95:D106").PasteSpecial(-4122) | Out-Null

foreach ($r in $newRows) {
    $row = $r.Row
    $time.Range("A$row").Value = $r.Task
    $time.Range("B$row").Value = $r.Start
    $time.Range("C$row").Value = $r.Stop
    $time.Range("D$row").Formula = "=IFERROR(MINUTE(C$row-B$row)/60+HOUR(C$row-B$row)+SECOND(C$row-B$row)/3600,0)"
}

$time.Columns.Item(1).AutoFit() | Out-Null

# ---------------------------------------------------------------------
# 2. Report sheet: task list text + new "Modified" dates for the two
#    in-flight tasks, now that the math-library research resolved them.
# ---------------------------------------------------------------------
$report.Range("A8").Value = "Math Unit Test and code"
$report.Range("E8").Value = 41730
$report.Range("E9").Value = 41730

# ---------------------------------------------------------------------
# 3. Report sheet: this week's narrative write-up.
# ---------------------------------------------------------------------
$report.Range("A25").Value = "Scripture`n"" 15 For behold, we are in bondage to the [senior project] and are taxed with a tax which is grievous to be borne..."" - Mosiah 7:15`n""Even the youths shall faint and be weary, and the young men shall utterly fall: "" Isaiah 40:30`n"
$report.Range("A22").Value = "On schedule for hours"
$report.Range("A15").Value = "Problems: Baby not letting me or my wife sleep as much as we want. When out of town this weekend. Some more research may be required. Due to finding great math library may not spend lots of time on math unit tests. Hopefully other areas will absorb more time."
$report.Range("A13").Value = "Accomplishments: First of many unit tests have been written, good start on Sequential Ray tracer, found great math library that should do most everything I need done. Managed not to get behind on hours. Design mostly finished."
$report.Range("A23").Value = "Will try to catch up on some time lost this week so I can be ahead of hours again."

# ---------------------------------------------------------------------
# 4. Workbook: drop the stale log_2 defined name / external query link,
#    and leave "Time" as the active, selected sheet (this week's focus).
# ---------------------------------------------------------------------
$wb.Names.Item("Time!log_2").Delete()

$report.Range("C5").Select() | Out-Null
$time.Activate() | Out-Null
$time.Range("A2").Select() | Out-Null

$wb.Application.Calculate()
